# Apply the "version final sin errores" edit:
#  - Metadata sheet: Version value 0.4.0 -> 0.7.0
#  - Metadata sheet: remove the "Jurisdiction" / "Chile" row entirely,
#    shifting the rows below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update the Version value (row 3, column B) before the row shift.
$ws.Range("B3").Value = "0.7.0"

# Delete the entire "Jurisdiction" row (row 11) and shift the remaining
# rows up, just like Excel's own Delete Row command would.
$ws.Rows.Item(11).Delete()
